# Auto-generated Excel COM-interop script to apply Ramuh_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 522.7692
$ws.Range("I107").Value = 421.1111
$ws.Range("J107").Value = 751.5
$ws.Range("K107").Value = 421.1111
$ws.Range("L107").Value = 751.5
$ws.Range("M107").Value = 1498.8889
$ws.Range("N107").Value = -4591.5

$ws.Range("H127").Value = 1167.8572
$ws.Range("I127").Value = 813.63635
$ws.Range("J127").Value = 2466.6667
$ws.Range("K127").Value = 2440.90905
$ws.Range("L127").Value = 7400.000100000001
$ws.Range("M127").Value = 2519.09095
$ws.Range("N127").Value = -17320.0001

$ws.Range("H138").Value = 2003.6765
$ws.Range("I138").Value = 803.9756
$ws.Range("J138").Value = 3825.4443
$ws.Range("K138").Value = 2411.9268
$ws.Range("L138").Value = 11476.3329
$ws.Range("M138").Value = 2728.0732
$ws.Range("N138").Value = -21756.3329

$ws.Range("H141").Value = 3148.5557
$ws.Range("I141").Value = 1048.3024
$ws.Range("J141").Value = 11358.637
$ws.Range("K141").Value = 3144.9072
$ws.Range("L141").Value = 34075.911
$ws.Range("M141").Value = 2035.0928
$ws.Range("N141").Value = -44435.911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 36000
$ws.Range("J107").Value = 36000
$ws.Range("L107").Value = 36000
$ws.Range("N107").Value = -43680

$ws.Range("H109").Value = 28450
$ws.Range("J109").Value = 28450
$ws.Range("L109").Value = 28450
$ws.Range("N109").Value = -31224

$ws.Range("H132").Value = 5596.7046
$ws.Range("I132").Value = 3861.0571
$ws.Range("J132").Value = 12346.444
$ws.Range("K132").Value = 11583.1713
$ws.Range("L132").Value = 37039.33199999999
$ws.Range("M132").Value = -9053.1713
$ws.Range("N132").Value = -42099.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 981.8333
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 981.8333
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 981.8333
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -3977.8333

$ws.Range("H107").Value = 1681.4584
$ws.Range("I107").Value = 801.8461
$ws.Range("J107").Value = 2721
$ws.Range("K107").Value = 801.8461
$ws.Range("L107").Value = 2721
$ws.Range("M107").Value = 1118.1539
$ws.Range("N107").Value = -6561

$ws.Range("H108").Value = 44561.332
$ws.Range("J108").Value = 44561.332
$ws.Range("L108").Value = 44561.332
$ws.Range("N108").Value = -52241.332

$ws.Range("H134").Value = 916.23914
$ws.Range("I134").Value = 802.9474
$ws.Range("J134").Value = 1454.375
$ws.Range("K134").Value = 2408.8422
$ws.Range("L134").Value = 4363.125
$ws.Range("M134").Value = 126.1578
$ws.Range("N134").Value = -9433.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1020.56714
$ws.Range("I58").Value = 758.56525
$ws.Range("J58").Value = 1594.4762
$ws.Range("K58").Value = 758.56525
$ws.Range("L58").Value = 1594.4762
$ws.Range("M58").Value = -555.56525
$ws.Range("N58").Value = -2000.4762

$ws.Range("H99").Value = 2821.4285
$ws.Range("I99").Value = 1684
$ws.Range("J99").Value = 4338
$ws.Range("K99").Value = 1684
$ws.Range("L99").Value = 4338
$ws.Range("M99").Value = -186
$ws.Range("N99").Value = -7334

$ws.Range("H107").Value = 938.1579
$ws.Range("I107").Value = 884.0833
$ws.Range("J107").Value = 1030.8572
$ws.Range("K107").Value = 884.0833
$ws.Range("L107").Value = 1030.8572
$ws.Range("M107").Value = 1035.9167
$ws.Range("N107").Value = -4870.8572

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H120").Value = 16400
$ws.Range("J120").Value = 16400
$ws.Range("L120").Value = 16400
$ws.Range("N120").Value = -23658

$ws.Range("H126").Value = 2821.4285
$ws.Range("I126").Value = 1684
$ws.Range("J126").Value = 4338
$ws.Range("K126").Value = 5052
$ws.Range("L126").Value = 13014
$ws.Range("M126").Value = -2582
$ws.Range("N126").Value = -17954

$ws.Range("H134").Value = 3931.5
$ws.Range("I134").Value = 4655.4546
$ws.Range("J134").Value = 1759.6364
$ws.Range("K134").Value = 13966.3638
$ws.Range("L134").Value = 5278.9092
$ws.Range("M134").Value = -11431.3638
$ws.Range("N134").Value = -10348.9092

$ws.Range("H136").Value = 1020.56714
$ws.Range("I136").Value = 758.56525
$ws.Range("J136").Value = 1594.4762
$ws.Range("K136").Value = 2275.69575
$ws.Range("L136").Value = 4783.4286
$ws.Range("M136").Value = 274.3042500000001
$ws.Range("N136").Value = -9883.428599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 476801.06
$ws.Range("I5").Value = 301.08334
$ws.Range("J5").Value = 1112134.4
$ws.Range("K5").Value = 903.2500200000001
$ws.Range("L5").Value = 3336403.2
$ws.Range("M5").Value = -791.2500200000001
$ws.Range("N5").Value = -3336627.2

$ws.Range("H107").Value = 273.9565
$ws.Range("I107").Value = 228
$ws.Range("J107").Value = 286.72223
$ws.Range("K107").Value = 684
$ws.Range("L107").Value = 860.16669
$ws.Range("M107").Value = 1236
$ws.Range("N107").Value = -4700.16669

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()

$ws.Range("H131").Value = 798.4353
$ws.Range("I131").Value = 508.87097
$ws.Range("J131").Value = 964.6667
$ws.Range("K131").Value = 1526.61291
$ws.Range("L131").Value = 2894.0001
$ws.Range("M131").Value = 3513.38709
$ws.Range("N131").Value = -12974.0001

$ws.Range("H135").Value = 476801.06
$ws.Range("I135").Value = 301.08334
$ws.Range("J135").Value = 1112134.4
$ws.Range("K135").Value = 2709.75006
$ws.Range("L135").Value = 10009209.6
$ws.Range("M135").Value = -174.7500600000003
$ws.Range("N135").Value = -10014279.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3157.2927
$ws.Range("I102").Value = 3078.375
$ws.Range("J102").Value = 3437.889
$ws.Range("K102").Value = 3078.375
$ws.Range("L102").Value = 3437.889
$ws.Range("M102").Value = -1456.375
$ws.Range("N102").Value = -6681.889

$ws.Range("H107").Value = 152.81482
$ws.Range("I107").Value = 146.63637
$ws.Range("J107").Value = 180
$ws.Range("K107").Value = 146.63637
$ws.Range("L107").Value = 180
$ws.Range("M107").Value = 1773.36363
$ws.Range("N107").Value = -4020

$ws.Range("H108").Value = 35242
$ws.Range("J108").Value = 35242
$ws.Range("L108").Value = 35242
$ws.Range("N108").Value = -42922

$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365

$ws.Range("H126").Value = 1814.2858
$ws.Range("I126").Value = 1480
$ws.Range("J126").Value = 2650
$ws.Range("K126").Value = 4440
$ws.Range("L126").Value = 7950
$ws.Range("M126").Value = -1970
$ws.Range("N126").Value = -12890

$ws.Range("H135").Value = 29930.525
$ws.Range("J135").Value = 29930.525
$ws.Range("L135").Value = 29930.525
$ws.Range("N135").Value = -40070.525

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2124.2666
$ws.Range("I93").Value = 2081.6
$ws.Range("J93").Value = 2145.6
$ws.Range("K93").Value = 2081.6
$ws.Range("L93").Value = 2145.6
$ws.Range("M93").Value = -833.5999999999999
$ws.Range("N93").Value = -4641.6

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H132").Value = 2701.2856
$ws.Range("I132").Value = 2648.4375
$ws.Range("J132").Value = 2961.4614
$ws.Range("K132").Value = 7945.3125
$ws.Range("L132").Value = 8884.3842
$ws.Range("M132").Value = -5415.3125
$ws.Range("N132").Value = -13944.3842

$ws.Range("H136").Value = 3085.4443
$ws.Range("I136").Value = 1168.9574
$ws.Range("J136").Value = 15953.286
$ws.Range("K136").Value = 3506.8722
$ws.Range("L136").Value = 47859.858
$ws.Range("M136").Value = -956.8721999999998
$ws.Range("N136").Value = -52959.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1460
$ws.Range("I96").Value = 900
$ws.Range("J96").Value = 1600
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 1600
$ws.Range("M96").Value = 473
$ws.Range("N96").Value = -4346

$ws.Range("H107").Value = 285.7
$ws.Range("I107").Value = 275.875
$ws.Range("J107").Value = 325
$ws.Range("K107").Value = 827.625
$ws.Range("L107").Value = 975
$ws.Range("M107").Value = 1092.375
$ws.Range("N107").Value = -4815

$ws.Range("H113").Value = 143345.42
$ws.Range("I113").Value = 125260.875
$ws.Range("K113").Value = 375782.625
$ws.Range("M113").Value = -373612.625

$ws.Range("H122").Value = 5129725.5
$ws.Range("I122").Value = 6452876
$ws.Range("J122").Value = 2517.5
$ws.Range("K122").Value = 19358628
$ws.Range("L122").Value = 7552.5
$ws.Range("M122").Value = -19356178
$ws.Range("N122").Value = -12452.5

$ws.Range("H132").Value = 1904.1023
$ws.Range("I132").Value = 2463.8928
$ws.Range("J132").Value = 924.46875
$ws.Range("K132").Value = 7391.678400000001
$ws.Range("L132").Value = 2773.40625
$ws.Range("M132").Value = -4861.678400000001
$ws.Range("N132").Value = -7833.40625

$ws.Range("H136").Value = 1065.5205
$ws.Range("I136").Value = 456.56818
$ws.Range("J136").Value = 1989.4482
$ws.Range("K136").Value = 1369.70454
$ws.Range("L136").Value = 5968.3446
$ws.Range("M136").Value = 1180.29546
$ws.Range("N136").Value = -11068.3446
